$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the value in E2 (was 04104008319, now 04104012998)
$ws.Range("E2").Value = "'04104012998"

# Update the active selection to E2 (was H11)
$ws.Range("E2").Select()
